$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "57.146.11"
$ws.Range("E2").Value = "  +5.25%  "
$ws.Range("D3").Value = "2.336.96"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.68"
$ws.Range("E5").Value = "  +4.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.31"
$ws.Range("E6").Value = "  +4.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.540"
$ws.Range("E8").Value = "  +2.69%  "
$ws.Range("D9").Value = "2.368.40"
$ws.Range("E9").Value = "  +4.01%  "
$ws.Range("E10").Value = "  +9.25%  "
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("E12").Value = "  +6.40%  "
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("E14").Value = "  +3.55%  "
$ws.Range("D15").Value = "2.781.25"
$ws.Range("E15").Value = "  +4.08%  "
$ws.Range("D16").Value = "57.184.52"
$ws.Range("E16").Value = "  +5.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  +5.21%  "
$ws.Range("D18").Value = "2.341.86"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.63"
$ws.Range("E19").Value = "  +3.38%  "
$ws.Range("E20").Value = "  +3.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.89"
$ws.Range("E21").Value = "  +6.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.67"
$ws.Range("E22").Value = "  +5.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.56"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.160"
$ws.Range("E25").Value = "  +7.44%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.990"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.82"
$ws.Range("E27").Value = "  +6.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.16"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").Value = "0.0₃0746"
$ws.Range("E29").Value = "  +6.22%  "
$ws.Range("E30").Value = "  +11.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.34"
$ws.Range("E31").Value = "  +5.45%  "
$ws.Range("E32").Value = "  +5.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.45"
$ws.Range("E33").Value = "  +3.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.958"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.993"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  +5.66%  "
$ws.Range("E38").Value = "  +8.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.53"
$ws.Range("E39").Value = "  +8.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.69"
$ws.Range("E40").Value = "  +4.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.385"
$ws.Range("E41").Value = "  +2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.42"
$ws.Range("E42").Value = "  +13.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.62"
$ws.Range("E43").Value = "  +6.89%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "280.81"
$ws.Range("E44").Value = "  +14.45%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.22"
$ws.Range("E45").Value = "  +8.08%  "
$ws.Range("E46").Value = "  +4.08%  "
$ws.Range("E47").Value = "  +3.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.568"
$ws.Range("E48").Value = "  +3.93%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0217"
$ws.Range("E49").Value = "  +6.04%  "
$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.384"
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.07"
$ws.Range("E51").Value = "  +4.99%  "
